# Generate Report for Handback
# Updates the localization-status report after a de-de handback:
#  - Overview/zh-cn/de-de "Status" columns move from "Ready for handoff"
#    to "Handed back: in sync with en-US"
#  - zh-cn/de-de rows gain a populated "Latest Target File" / "Latest
#    Handback File" / "Latest Handback DateTime" for the handed-back file
#  - Column widths are widened to fit the new values

$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ab3f7d77544159f3a73526aaf848146fc22f8e2f/e2e/ecd86253-0024-480b-a560-9d81cf9ba420.md"
$targetDisplay = "ecd86253-0024-480b-a560-9d81cf9ba420.md"
$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsZhCn.Range("I2").Value = $targetDisplay
$wsZhCn.Range("J2").Value = "ecd86253-0024-480b-a560-9d81cf9ba420.35dd873a71607bc66db71ef1ca7c83543327024c.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-06 03:11:17"

$wsZhCn.Range("I3").Value = $targetDisplay
$wsZhCn.Range("J3").Value = "ecd86253-0024-480b-a560-9d81cf9ba420.35dd873a71607bc66db71ef1ca7c83543327024c.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-06 03:11:17"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $targetUrl, "", "", $targetDisplay) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $targetUrl, "", "", $targetDisplay) | Out-Null

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

$wsDeDe.Range("I2").Value = $targetDisplay
$wsDeDe.Range("J2").Value = "ecd86253-0024-480b-a560-9d81cf9ba420.35dd873a71607bc66db71ef1ca7c83543327024c.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-06 03:11:24"

$wsDeDe.Range("I3").Value = $targetDisplay
$wsDeDe.Range("J3").Value = "ecd86253-0024-480b-a560-9d81cf9ba420.35dd873a71607bc66db71ef1ca7c83543327024c.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-06 03:11:24"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $targetUrl, "", "", $targetDisplay) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $targetUrl, "", "", $targetDisplay) | Out-Null

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664
